# Update the "data" sheet time_taken (column F) timestamps for rows 2..115,
# then add a new "metadata" sheet (placed after "data") describing the
# panelapp query that produced this export.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$newTimes = @(
    "2021-10-05 14:19:19.719261",
    "2021-10-05 14:19:19.719269",
    "2021-10-05 14:19:19.719272",
    "2021-10-05 14:19:19.719275",
    "2021-10-05 14:19:19.719278",
    "2021-10-05 14:19:19.719280",
    "2021-10-05 14:19:19.719283",
    "2021-10-05 14:19:19.719285",
    "2021-10-05 14:19:19.719288",
    "2021-10-05 14:19:19.719291",
    "2021-10-05 14:19:19.719293",
    "2021-10-05 14:19:19.719295",
    "2021-10-05 14:19:19.719298",
    "2021-10-05 14:19:19.719300",
    "2021-10-05 14:19:19.719303",
    "2021-10-05 14:19:19.719305",
    "2021-10-05 14:19:19.719308",
    "2021-10-05 14:19:19.719311",
    "2021-10-05 14:19:19.719313",
    "2021-10-05 14:19:19.719316",
    "2021-10-05 14:19:19.719318",
    "2021-10-05 14:19:19.719321",
    "2021-10-05 14:19:19.719323",
    "2021-10-05 14:19:19.719326",
    "2021-10-05 14:19:19.719329",
    "2021-10-05 14:19:19.719331",
    "2021-10-05 14:19:19.719334",
    "2021-10-05 14:19:19.719336",
    "2021-10-05 14:19:19.719339",
    "2021-10-05 14:19:19.719341",
    "2021-10-05 14:19:19.719344",
    "2021-10-05 14:19:19.719346",
    "2021-10-05 14:19:19.719349",
    "2021-10-05 14:19:19.719351",
    "2021-10-05 14:19:19.719354",
    "2021-10-05 14:19:19.719356",
    "2021-10-05 14:19:19.719359",
    "2021-10-05 14:19:19.719361",
    "2021-10-05 14:19:19.719364",
    "2021-10-05 14:19:19.719366",
    "2021-10-05 14:19:19.719369",
    "2021-10-05 14:19:19.719372",
    "2021-10-05 14:19:19.719374",
    "2021-10-05 14:19:19.719377",
    "2021-10-05 14:19:19.719379",
    "2021-10-05 14:19:19.719382",
    "2021-10-05 14:19:19.719384",
    "2021-10-05 14:19:19.719386",
    "2021-10-05 14:19:19.719389",
    "2021-10-05 14:19:19.719391",
    "2021-10-05 14:19:19.719394",
    "2021-10-05 14:19:19.719396",
    "2021-10-05 14:19:19.719399",
    "2021-10-05 14:19:19.719402",
    "2021-10-05 14:19:19.719404",
    "2021-10-05 14:19:19.719407",
    "2021-10-05 14:19:19.719409",
    "2021-10-05 14:19:19.719412",
    "2021-10-05 14:19:19.719414",
    "2021-10-05 14:19:19.719417",
    "2021-10-05 14:19:19.719419",
    "2021-10-05 14:19:19.719421",
    "2021-10-05 14:19:19.719424",
    "2021-10-05 14:19:19.719426",
    "2021-10-05 14:19:19.719430",
    "2021-10-05 14:19:19.719433",
    "2021-10-05 14:19:19.719436",
    "2021-10-05 14:19:19.719438",
    "2021-10-05 14:19:19.719441",
    "2021-10-05 14:19:19.719443",
    "2021-10-05 14:19:19.719445",
    "2021-10-05 14:19:19.719448",
    "2021-10-05 14:19:19.719450",
    "2021-10-05 14:19:19.719453",
    "2021-10-05 14:19:19.719455",
    "2021-10-05 14:19:19.719458",
    "2021-10-05 14:19:19.719462",
    "2021-10-05 14:19:19.719465",
    "2021-10-05 14:19:19.719468",
    "2021-10-05 14:19:19.719470",
    "2021-10-05 14:19:19.719473",
    "2021-10-05 14:19:19.719475",
    "2021-10-05 14:19:19.719478",
    "2021-10-05 14:19:19.719480",
    "2021-10-05 14:19:19.719483",
    "2021-10-05 14:19:19.719485",
    "2021-10-05 14:19:19.719488",
    "2021-10-05 14:19:19.719490",
    "2021-10-05 14:19:19.719493",
    "2021-10-05 14:19:19.719495",
    "2021-10-05 14:19:19.719497",
    "2021-10-05 14:19:19.719500",
    "2021-10-05 14:19:19.719503",
    "2021-10-05 14:19:19.719506",
    "2021-10-05 14:19:19.719509",
    "2021-10-05 14:19:19.719511",
    "2021-10-05 14:19:19.719514",
    "2021-10-05 14:19:19.719516",
    "2021-10-05 14:19:19.719519",
    "2021-10-05 14:19:19.719521",
    "2021-10-05 14:19:19.719523",
    "2021-10-05 14:19:19.719526",
    "2021-10-05 14:19:19.719528",
    "2021-10-05 14:19:19.719531",
    "2021-10-05 14:19:19.719533",
    "2021-10-05 14:19:19.719536",
    "2021-10-05 14:19:19.719539",
    "2021-10-05 14:19:19.719541",
    "2021-10-05 14:19:19.719546",
    "2021-10-05 14:19:19.719549",
    "2021-10-05 14:19:19.719552",
    "2021-10-05 14:19:19.719554",
    "2021-10-05 14:19:19.719557",
    "2021-10-05 14:19:19.719560"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 6).Value = $newTimes[$i]
}

# --- add the metadata sheet -------------------------------------------------
$new = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$new.Name = "metadata"

# Pull the header cell style (bold, bordered, centered) from the "data" sheet
# and apply it to the metadata header row + the A2 index-style cell, mirroring
# the layout of the "data" sheet.
$ws1.Range("B1").Copy()
$new.Range("B1:G1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$new.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$new.Range("B1").Value = "data_name"
$new.Range("C1").Value = "data_id"
$new.Range("D1").Value = "data_version"
$new.Range("E1").Value = "data_version_created"
$new.Range("F1").Value = "panel_query_time"
$new.Range("G1").Value = "panel_get_request"

$new.Range("A2").Value = 0
$new.Range("B2").Value = "Bleeding and platelet disorders"
$new.Range("C2").Value = 545

# data_version ("1.31") must stay text, not be coerced to the number 1.31
$new.Range("D2").NumberFormat = "@"
$new.Range("D2").Value = "1.31"
$new.Range("D2").ClearFormats()

$new.Range("E2").Value = "2021-08-02T07:42:26.362060Z"
$new.Range("F2").Value = "2021-10-05 14:19:19.715706"
$new.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/545/?format=json"

# Keep "data" as the active/selected sheet (matches activeTab=0 upstream).
$ws1.Activate()
$ws1.Range("A1").Select() | Out-Null
